$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Continue Running" config row first: label in A4, boolean TRUE in B4
$ws.Range("A4").Value = "Continue Running"
$ws.Range("B4").Value = $true

# Update the Arrival Station value from the placeholder "Adsfdsfa" to "Constanta"
$ws.Range("B3").Value = "Constanta"

# Update the selection to reflect the new active cell/row as in the diff
$ws.Range("A5:B5").Select()
